# Qualia sample template — Dec 2024 ops changes
# (ordercreation & orderpage & order form)
#
# Reorders the Typist/Client/Lob/Process/Product-Name header block,
# replaces the sample order rows with the new SIPL/Qu18 sample data,
# and adds a new "Tier" column (N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Capture the "emphasis" cell format (bordered + special font, no
# left edge trimmed off) from the existing Client column (E2/E3)
# *before* anything else changes, so we can stamp it onto the cells
# that need it at their new positions (C, D, G).
# -----------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("G3").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Borders.Item(7).LineStyle = -4142
$ws.Range("E3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Borders.Item(7).LineStyle = -4142

# Now that its look has been copied onward, put E2/E3 back to the
# plain bordered style shared by the rest of the data cells (M2's
# style, which is the same plain bordered look as B/F/H/I/J/K/L/M).
$ws.Range("M2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("M3").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Header row (row 1): columns E..J get reordered, N1 is new ("Tier")
# ---------------------------------------------------------------
$ws.Range("E1").Value = "Typist"
$ws.Range("F1").Value = "Typist QC"
$ws.Range("G1").Value = "Client"
$ws.Range("H1").Value = "Lob"
$ws.Range("I1").Value = "Process"
$ws.Range("J1").Value = "Product Name"

# New Tier header — copy the format from the neighboring header cell (M1)
# so it picks up the same header style (bold font + fill + border).
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "Tier"

# ---------------------------------------------------------------
# Row 2 (first sample order)
# ---------------------------------------------------------------
$ws.Range("B2").Value = "Qu18-001"
$ws.Range("C2").Value = "SIPL6118"
$ws.Range("D2").Value = "SIPL4167"
$ws.Range("E2").Value = "SIPL5317"
$ws.Range("F2").Value = "SIPL5317"
$ws.Range("G2").Value = "Qualia"
$ws.Range("H2").Value = "Title"
$ws.Range("I2").Value = "Search & Typing"
$ws.Range("J2").Value = "Current Owner Search"
$ws.Range("K2").Value = "AL"
$ws.Range("L2").Value = "Shelby"
$ws.Range("M2").Value = "WIP"

# ---------------------------------------------------------------
# Row 3 (second sample order)
# ---------------------------------------------------------------
$ws.Range("B3").Value = "Qu18-002"
$ws.Range("C3").Value = "SIPL5316"
$ws.Range("D3").Value = "SIPL5688"
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL0103"
$ws.Range("G3").Value = "Qualia"
$ws.Range("H3").Value = "Title"
$ws.Range("I3").Value = "Search & Typing"
$ws.Range("J3").Value = "Full Search"
$ws.Range("K3").Value = "FL"
$ws.Range("L3").Value = "Clay"
$ws.Range("M3").Value = "WIP"

# ---------------------------------------------------------------
# New Tier cells (N2/N3) — left empty, but still bordered like the
# rest of the data rows so the table outline stays closed.
# ---------------------------------------------------------------
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("N2").ClearContents()

$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").ClearContents()

# ---------------------------------------------------------------
# Column widths (best-fit adjustments that came along with the
# header/content changes)
# ---------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 19.5
$ws.Columns.Item(7).ColumnWidth = 5.3333333333333333
$ws.Columns.Item(8).ColumnWidth = 3.5
$ws.Columns.Item(9).ColumnWidth = 12.8333333333333333
$ws.Columns.Item(10).ColumnWidth = 18.1666666666666667
$ws.Columns.Item(14).ColumnWidth = 10.6666666666666667

# ---------------------------------------------------------------
# Selection / active cell
# ---------------------------------------------------------------
$ws.Range("E5").Select()
